{"js": "// Apply the diversion-program wording change, fines/jail-days table\n// updates, the \"court costs\" sentence addition, and the small run\n// merge in the \"proof of financial responsibility\" paragraph.\n\nconst body = context.document.body;\n\n// --- 1) \"Prosecutor Diversion Program\" -> \"Marijuana Diversion Program\" ---\nlet results = body.search(\"Prosecutor Diversion Program\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Marijuana Diversion Program\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 2) Both remaining lower-case \"diversion program\" -> \"Marijuana Diversion Program\" ---\nresults = body.search(\"diversion program\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nfor (const r of results.items) {\n  r.insertText(\"Marijuana Diversion Program\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// --- 3) Table values: Fine Amount $50 -> $0, Fines Suspended $25 -> $0,\n//        Jail Days 10 -> None, Jail Days Suspended 7 -> None ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nasync function replaceInCell(cell, oldText, newText) {\n  const found = cell.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"text\");\n  await context.sync();\n  found.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// Row index 5 = \"Fine Amount\"\nlet cells = rows.items[5].cells;\ncells.load(\"items\");\nawait context.sync();\nawait replaceInCell(cells.items[1], \"$ 50\", \"$ 0\");\n\n// Row index 6 = \"Fines Suspended\"\ncells = rows.items[6].cells;\ncells.load(\"items\");\nawait context.sync();\nawait replaceInCell(cells.items[1], \"$ 25\", \"$ 0\");\n\n// Row index 7 = \"Jail Days\"\ncells = rows.items[7].cells;\ncells.load(\"items\");\nawait context.sync();\nawait replaceInCell(cells.items[1], \"10\", \"None\");\n\n// Row index 8 = \"Jail Days Suspended\"\ncells = rows.items[8].cells;\ncells.load(\"items\");\nawait context.sync();\nawait replaceInCell(cells.items[1], \"7\", \"None\");\n\n// --- 4) \"...in this case\" -> \"...in this case, but are not due if diversion is successfully completed\" ---\n//        (search text excludes the trailing period on purpose: the period belongs to a\n//        separate, non-underlined run that must stay untouched / not become underlined)\nresults = body.search(\"in this case\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\n  \"in this case, but are not due if diversion is successfully completed\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// --- 5) Remove the mid-sentence page-break split: merge into one continuous run/sentence ---\nresults = body.search(\"The Defendant showed the Court proof of responsibility during the proceeding.\", {\n  matchCase: true,\n  matchWholeWord: false,\n});\nresults.load(\"text\");\nawait context.sync();\nresults.items[0].insertText(\n  \"The Defendant showed the Court proof of responsibility during the proceeding.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "# Apply the diversion-program wording change, fines/jail-days table\n# updates, the \"court costs\" sentence addition, and the small run\n# merge in the \"proof of financial responsibility\" paragraph.\n\n$d = $word.ActiveDocument\n\nfunction Replace-All($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Execute(\n        $findText,    # FindText\n        $true,        # MatchCase\n        $false,       # MatchWholeWord\n        $false,       # MatchWildcards\n        $false,       # MatchSoundsLike\n        $false,       # MatchAllWordForms\n        $true,        # Forward\n        1,            # Wrap (wdFindContinue)\n        $false,       # Format\n        $replaceText, # ReplaceWith\n        2             # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n\n# --- 1) \"Prosecutor Diversion Program\" -> \"Marijuana Diversion Program\" ---\nReplace-All \"Prosecutor Diversion Program\" \"Marijuana Diversion Program\"\n\n# --- 2) Both remaining lower-case \"diversion program\" -> \"Marijuana Diversion Program\" ---\nReplace-All \"diversion program\" \"Marijuana Diversion Program\"\n\n# --- 3) Table values: Fine Amount $50 -> $0, Fines Suspended $25 -> $0,\n#        Jail Days 10 -> None, Jail Days Suspended 7 -> None ---\n$t = $d.Tables.Item(1)\n$t.Cell(6, 2).Range.Text = \"$ 0\"\n$t.Cell(7, 2).Range.Text = \"$ 0\"\n$t.Cell(8, 2).Range.Text = \"None\"\n$t.Cell(9, 2).Range.Text = \"None\"\n\n# --- 4) \"...in this case\" -> \"...in this case, but are not due if diversion is successfully completed\" ---\nReplace-All \"in this case\" \"in this case, but are not due if diversion is successfully completed\"\n\n# --- 5) Remove the mid-sentence page-break split: merge into one continuous run/sentence ---\nReplace-All \"The Defendant showed the Court proof of responsibility during the proceeding.\" \"The Defendant showed the Court proof of responsibility during the proceeding.\"\n"}
